$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Re-run of the RAD Katalon test suite: refresh the "Date" (Execute timestamp)
# column for the previously recorded rows (2-54) to reflect the new run.
$dateUpdates = @{
    2 = "Tue Feb 11 20:07:37 EST 2025"
    3 = "Tue Feb 11 20:07:49 EST 2025"
    4 = "Tue Feb 11 20:08:01 EST 2025"
    5 = "Tue Feb 11 20:08:13 EST 2025"
    6 = "Tue Feb 11 20:08:25 EST 2025"
    7 = "Tue Feb 11 20:08:36 EST 2025"
    8 = "Tue Feb 11 20:08:48 EST 2025"
    9 = "Tue Feb 11 20:09:00 EST 2025"
    10 = "Tue Feb 11 20:09:12 EST 2025"
    11 = "Tue Feb 11 20:09:24 EST 2025"
    12 = "Tue Feb 11 20:09:36 EST 2025"
    13 = "Tue Feb 11 20:09:48 EST 2025"
    14 = "Tue Feb 11 20:10:00 EST 2025"
    15 = "Tue Feb 11 20:10:12 EST 2025"
    16 = "Tue Feb 11 20:10:25 EST 2025"
    17 = "Tue Feb 11 20:10:38 EST 2025"
    18 = "Tue Feb 11 20:10:50 EST 2025"
    19 = "Tue Feb 11 20:11:03 EST 2025"
    20 = "Tue Feb 11 20:11:15 EST 2025"
    21 = "Tue Feb 11 20:11:28 EST 2025"
    22 = "Tue Feb 11 20:11:40 EST 2025"
    23 = "Tue Feb 11 20:11:53 EST 2025"
    24 = "Tue Feb 11 20:12:06 EST 2025"
    25 = "Tue Feb 11 20:12:19 EST 2025"
    26 = "Tue Feb 11 20:12:31 EST 2025"
    27 = "Tue Feb 11 20:12:43 EST 2025"
    28 = "Tue Feb 11 20:12:56 EST 2025"
    29 = "Tue Feb 11 20:13:09 EST 2025"
    30 = "Tue Feb 11 20:13:22 EST 2025"
    31 = "Tue Feb 11 20:13:35 EST 2025"
    32 = "Tue Feb 11 20:13:48 EST 2025"
    33 = "Tue Feb 11 20:14:00 EST 2025"
    34 = "Tue Feb 11 20:14:13 EST 2025"
    35 = "Tue Feb 11 20:14:25 EST 2025"
    36 = "Tue Feb 11 20:14:37 EST 2025"
    37 = "Tue Feb 11 20:14:49 EST 2025"
    38 = "Tue Feb 11 20:15:01 EST 2025"
    39 = "Tue Feb 11 20:15:12 EST 2025"
    40 = "Tue Feb 11 20:15:24 EST 2025"
    41 = "Tue Feb 11 20:15:36 EST 2025"
    42 = "Tue Feb 11 20:15:48 EST 2025"
    43 = "Tue Feb 11 20:16:00 EST 2025"
    44 = "Tue Feb 11 20:16:12 EST 2025"
    45 = "Tue Feb 11 20:16:25 EST 2025"
    46 = "Tue Feb 11 20:16:38 EST 2025"
    47 = "Tue Feb 11 20:16:50 EST 2025"
    48 = "Tue Feb 11 20:17:03 EST 2025"
    49 = "Tue Feb 11 20:17:16 EST 2025"
    50 = "Tue Feb 11 20:17:28 EST 2025"
    51 = "Tue Feb 11 20:17:41 EST 2025"
    52 = "Tue Feb 11 20:17:54 EST 2025"
    53 = "Tue Feb 11 20:18:06 EST 2025"
    54 = "Tue Feb 11 20:18:19 EST 2025"
}
foreach ($row in $dateUpdates.Keys) {
    $ws.Cells.Item($row, 2).Value2 = $dateUpdates[$row]
}

# Append newly recorded test rows (55-59) captured during this run, covering
# the additional PaymentType/TaxType combinations that were exercised
# (Digital Advertising Gross Revenues / PTE Composite / IFTA Tax).
$newRows = @(
    [PSCustomObject]@{ Row = 55; Date = "Tue Feb 11 20:18:32 EST 2025"; PaymentType = "Existing Liability with Notice/Invoice Number"; TaxType = "Digital Advertising Gross Revenues" }
    [PSCustomObject]@{ Row = 56; Date = "Tue Feb 11 20:18:44 EST 2025"; PaymentType = "New Tax Return Amount Due"; TaxType = "Digital Advertising Gross Revenues" }
    [PSCustomObject]@{ Row = 57; Date = "Tue Feb 11 20:18:56 EST 2025"; PaymentType = "Existing Liability with Notice/Invoice Number"; TaxType = "PTE Composite" }
    [PSCustomObject]@{ Row = 58; Date = "Tue Feb 11 20:19:08 EST 2025"; PaymentType = "New Tax Return Amount Due"; TaxType = "IFTA Tax" }
    [PSCustomObject]@{ Row = 59; Date = "Tue Feb 11 20:19:21 EST 2025"; PaymentType = "New Tax Return Amount Due"; TaxType = "PTE Composite" }
)
foreach ($nr in $newRows) {
    $row = $nr.Row
    $ws.Cells.Item($row, 1).Value2 = "Pass"
    $ws.Cells.Item($row, 1).Style = "Normal"
    $ws.Cells.Item($row, 2).Value2 = $nr.Date
    $ws.Cells.Item($row, 2).Style = "Normal"
    $ws.Cells.Item($row, 3).Value2 = "Y"
    $ws.Cells.Item($row, 4).Value2 = $nr.PaymentType
    $ws.Cells.Item($row, 5).Value2 = $nr.TaxType
}

# Reflect the post-edit selection left by the author (column C across the newly
# touched/added rows).
$ws.Range("C46:C59").Select()
